$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Column indices (1-based) within each data row:
#  1=data 2=NumPedido 3=Cliente 4=Prazo 5=ValorPedido 6=frete
#  7=Referencia 8=%comissao 9=ValorComissao 10=Pagamento 11=Obs

$row = $t.Rows.Item(2)
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "27.539,62"

$row = $t.Rows.Item(3)
$row.Cells.Item(4).Range.Text = "30 a 90"
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "15.496,27"

$row = $t.Rows.Item(4)
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "17.554,22"

$row = $t.Rows.Item(5)
$row.Cells.Item(4).Range.Text = "15 a 45"
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "9.611,11"

$row = $t.Rows.Item(6)
$row.Cells.Item(4).Range.Text = "15 a 45"
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "17.019,02"

$row = $t.Rows.Item(7)
$row.Cells.Item(4).Range.Text = "30 a 90"
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "32.146,05"

$row = $t.Rows.Item(8)
$row.Cells.Item(4).Range.Text = "30 a 90"
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "12.444,52"

$row = $t.Rows.Item(9)
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "28.188,30"

$row = $t.Rows.Item(10)
$row.Cells.Item(4).Range.Text = "30 a 120"
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "19.481,41"

$row = $t.Rows.Item(11)
$row.Cells.Item(4).Range.Text = "30 a 120"
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "14.910,08"

$row = $t.Rows.Item(12)
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "15.309,34"

$row = $t.Rows.Item(13)
$row.Cells.Item(4).Range.Text = "30 a 120"
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "16.031,43"

$row = $t.Rows.Item(14)
$row.Cells.Item(4).Range.Text = "30 a 90"
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "26.175,60"

$row = $t.Rows.Item(15)
$row.Cells.Item(4).Range.Text = "30 a 90"
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "43.124,40"

$row = $t.Rows.Item(16)
$row.Cells.Item(4).Range.Text = "30 a 90"
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "16.544,00"

$row = $t.Rows.Item(17)
$row.Cells.Item(4).Range.Text = "30 a 90"
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "17.309,60"

$row = $t.Rows.Item(18)
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "8.683,41"

$row = $t.Rows.Item(19)
$row.Cells.Item(4).Range.Text = "30 a 60"
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "6.783,14"

$row = $t.Rows.Item(20)
$row.Cells.Item(4).Range.Text = "30 a 90"
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "55.474,50"

$row = $t.Rows.Item(21)
$row.Cells.Item(4).Range.Text = "30 a 90"
$row.Cells.Item(8).Range.Text = "5"
$row.Cells.Item(9).Range.Text = "13.536,64"
